# HCR 7: no fishing added
# Adds three new rows (51-53) describing the new "HCR7" (no fishing)
# scenario combinations, mirroring the layout/styling of the existing
# "INNvar" blocks (e.g. rows 30-32) but without the merged/commented
# H column used by the other blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (fill/font/alignment) of an existing "var" block
# (rows 30:32 -> A:G) down onto the new rows 51:53, so the new rows pick
# up the same alternating style pattern used elsewhere in the sheet.
$ws.Range("A30:G32").Copy() | Out-Null
$ws.Range("A51:G53").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Row 51 -> id 37, HCR7, REClow
$ws.Range("A51").Value = 37
$ws.Range("B51").Value = "var"
$ws.Range("C51").Value = "low"
$ws.Range("D51").Value = "none"
$ws.Range("E51").Value = "none"
$ws.Range("F51").Value = 7
$ws.Range("G51").Formula = '=CONCATENATE("ASS",E51,"_HCR",F51,"_REC",C51,"_INN",B51,"_OER",D51)'

# Row 52 -> id 38, HCR7, RECmed
$ws.Range("A52").Value = 38
$ws.Range("B52").Value = "var"
$ws.Range("C52").Value = "med"
$ws.Range("D52").Value = "none"
$ws.Range("E52").Value = "none"
$ws.Range("F52").Value = 7
$ws.Range("G52").Formula = '=CONCATENATE("ASS",E52,"_HCR",F52,"_REC",C52,"_INN",B52,"_OER",D52)'

# Row 53 -> id 39, HCR7, RECmix
$ws.Range("A53").Value = 39
$ws.Range("B53").Value = "var"
$ws.Range("C53").Value = "mix"
$ws.Range("D53").Value = "none"
$ws.Range("E53").Value = "none"
$ws.Range("F53").Value = 7
$ws.Range("G53").Formula = '=CONCATENATE("ASS",E53,"_HCR",F53,"_REC",C53,"_INN",B53,"_OER",D53)'

# Match the author's final on-screen selection.
$ws.Range("G57").Select()
